$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting existing D:K data to E:L
$ws.Columns("D").Insert()

# Copy number formats/styles from column E (the old column D, now shifted) into the
# newly inserted (blank) column D so the new cells inherit the correct formatting
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the newly inserted column D with the new period's financial data
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 207500
$ws.Range("D9").Value = 96100
$ws.Range("D10").Value = 111400
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 99100
$ws.Range("D18").Value = 108400
$ws.Range("D20").Value = -69300
$ws.Range("D21").Value = 39100
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 39100
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 39100
$ws.Range("D27").Value = 39100
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 69300
$ws.Range("D33").Value = 39100
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 39100
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 87200
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 21300
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 100
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 1972200
$ws.Range("D48").Value = 0
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 4000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2084700
$ws.Range("D57").Value = 13900
$ws.Range("D58").Value = 446000
$ws.Range("D59").Value = 46900
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 514600
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1021500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -110100
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1063200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 39100
$ws.Range("D83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 28800
$ws.Range("D91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 0
$ws.Range("D96").Value = -93300
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 26400
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 55100
